# Fixed retrieve a patient's appointments for a patient where no appointments
# should be present: add a new "patientNoAppointments" prerequisite row to the
# Patients sheet, and make that sheet the active tab again (was left on
# "Practitioner" previously).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patients")

# Duplicate the last prerequisite row (22) into a new row 23, preserving all
# per-column formatting (borders/fills alternate across the row).
$ws.Range("B22:P22").Copy($ws.Range("B23:P23"))

# New prerequisite: a patient that must have no appointments at all.
$ws.Range("B23").Value = "patientNoAppointments"
$ws.Range("P23").Value = "Patient should not have any appointments"

# Match the taller row used for this new entry.
$ws.Rows.Item(23).RowHeight = 23.3

# Column B needs to be wider to fit the new, longer identifier text.
$ws.Columns.Item(2).ColumnWidth = 20

# The workbook had been left with "Practitioner" as the active/selected tab;
# switch back to "Patients" and leave the cursor on the freshly added row.
$ws.Activate()
$ws.Range("A23").Select()
